$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 167, shifting rows 167:177 down to 168:178
$ws.Rows.Item(167).Insert()

# Fill in the new row 167 data (mirrors the pattern of the surrounding "Jengibre" rows)
$ws.Cells.Item(167, 1).Value = 8
$ws.Cells.Item(167, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(167, 3).Value = "Coquimbo"
$ws.Cells.Item(167, 4).Value = 45166
$ws.Cells.Item(167, 4).NumberFormat = $ws.Cells.Item(168, 4).NumberFormat
$ws.Cells.Item(167, 5).Value = 4
$ws.Cells.Item(167, 6).Value = 100114007
$ws.Cells.Item(167, 7).Value = "Jengibre"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 400
$ws.Cells.Item(167, 11).Value = 16000
$ws.Cells.Item(167, 12).Value = 17000
$ws.Cells.Item(167, 13).Value = 16500
$ws.Cells.Item(167, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(167, 15).Value = "Perú"
$ws.Cells.Item(167, 16).Value = 1269
$ws.Cells.Item(167, 17).Value = 13
$ws.Cells.Item(167, 18).Value = "Hortaliza"
